$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 82

# Column A holds a value that looks like a date ("2025-10-16"). A plain
# Range.Value assignment lets Excel auto-detect it as a real date and
# reformat/restyle the cell. To keep it as plain text (matching the rest
# of the column, which stores these as inline/shared strings) without
# leaving a lingering custom cell style behind, stage the text in a
# scratch cell that is explicitly formatted as Text, copy/paste its
# value into place, then clean the scratch cell back up.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "2025-10-16"
$scratch.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)
$scratch.Clear()

$ws.Range("B" + $newRow).Value = "21:21:45"
$ws.Range("C" + $newRow).Value = "1.00 EUR = 1,703.2338"
